# "Add files via upload" -- the author re-uploaded/re-saved the workbook in
# Excel. The only deliberate content change captured by the diff is the
# worksheet being renamed from its default import name to "Sales"; the rest
# of the diff (new xr/xr2/xr3 revision namespaces, a regenerated
# xr:revisionPtr documentId/uid, and the workbookView's window geometry) is
# metadata Excel stamps on every save and isn't something the workbook
# content model exposes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (was "AdventureWorks_Sales_2017 (1)").
$ws.Name = "Sales"
